$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (styles) from the last existing row (550) down to the
# new rows (551:560) so the new rows inherit the same cell styles
# (date format on A, grey/white fonts on B-H, centered/merged style on G).
$ws.Range("A550:I550").Copy()
$ws.Range("A551:I560").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data for the 10 new rows (all dated 45973 = 2025-11-12)
$data = @(
    @{ Row = 551; Name = "Amir Etien";       C = 70; D = 9; E = 9; F = 0; G = $null;                     H = 8 },
    @{ Row = 552; Name = "Yoan Zouma";       C = 70; D = 7; E = 9; F = 5; G = "Ischio";                  H = 4 },
    @{ Row = 553; Name = "Yoann Martelat";   C = 70; D = 7; E = 7; F = 5; G = "Genou / Quadri";          H = 5 },
    @{ Row = 554; Name = "Ilyes Boughanmi";  C = 70; D = 8; E = 8; F = 2; G = "Genou";                   H = 0 },
    @{ Row = 555; Name = "Omar Benyounes";   C = 70; D = 7; E = 7; F = 2; G = "Cheville";                H = 6 },
    @{ Row = 556; Name = "Naim Ighbane";     C = 70; D = 7; E = 7; F = 7; G = "Cheville/adduct/genou";   H = 7 },
    @{ Row = 557; Name = "Karim Belmahi";    C = 70; D = 7; E = 8; F = 0; G = $null;                     H = 10 },
    @{ Row = 558; Name = "Malik Boussaid";   C = 70; D = 3; E = 0; F = 0; G = $null;                     H = 10 },
    @{ Row = 559; Name = "Emmanuel Valey";   C = 70; D = 7; E = 6; F = 4; G = "Ischio";                  H = 2 },
    @{ Row = 560; Name = "Sofiane Belle";    C = 70; D = 7; E = 6; F = 0; G = $null;                     H = 6 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = 45973
    $ws.Range("B$r").Value = $item.Name
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    if ($item.G) {
        $ws.Range("G$r").Value = $item.G
    }
    $ws.Range("H$r").Value = $item.H
    $ws.Range("I$r").Formula = "=C$r*D$r"
}

# Refresh the view selection to mirror the edited workbook
$ws.Range("L552").Select()

Write-Host "Applied wellness data update"
